# "hand routing the 6 layer stackup." — add a new "6-layer" worksheet
# (a trimmed-down copy of the "8-layer" stackup sheet) right after the
# existing "8-layer" sheet, and make it the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("8-layer")

# Insert the new sheet right after "8-layer" so it lands as tab #2.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "6-layer"

# ---- header row -----------------------------------------------------
$ws2.Range("B1").Value = "Copper Layer"
$ws2.Range("C1").Value = "Dielectric Layer"
$ws2.Range("D1").Value = "Material"
$ws2.Range("E1").Value = "Er"
$ws2.Range("F1").Value = "Thickness"
$ws2.Range("G1").Value = "Desired Z0"
$ws2.Range("H1").Value = "Trace Width"
$ws2.Range("I1").Value = "Calculated Z0"
$ws2.Range("A1:I1").HorizontalAlignment = -4108

# ---- stackup rows -----------------------------------------------------
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "top component"
$ws2.Range("F2").Value = 1
$ws2.Range("H2").Value = 8
$ws2.Range("I2").Value = 54

$ws2.Range("C3").Value = "prepreg"
$ws2.Range("D3").Value = "Grace GA-170LL"
$ws2.Range("E3").Value = 4.7
$ws2.Range("F3").Value = 6

$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = "ground plane"
$ws2.Range("F4").Value = 1

$ws2.Range("C5").Value = "laminate"
$ws2.Range("D5").Value = "Grace GA-170LL"
$ws2.Range("E5").Value = 4.7
$ws2.Range("F5").Value = 6

$ws2.Range("A6").Value = 4
$ws2.Range("B6").Value = "power plane"
$ws2.Range("F6").Value = 1

$ws2.Range("C7").Value = "prepreg"
$ws2.Range("D7").Value = "Grace GA-170LL"
$ws2.Range("E7").Value = 4.7
$ws2.Range("F7").Value = 14

$ws2.Range("A8").Value = 5
$ws2.Range("B8").Value = "power plane"
$ws2.Range("F8").Value = 1

$ws2.Range("C9").Value = "laminate"
$ws2.Range("D9").Value = "Grace GA-170LL"
$ws2.Range("E9").Value = 4.7
$ws2.Range("F9").Value = 6

$ws2.Range("A10").Value = 6
$ws2.Range("B10").Value = "ground plane"
$ws2.Range("F10").Value = 1

$ws2.Range("C11").Value = "prepreg"
$ws2.Range("D11").Value = "Grace GA-170LL"
$ws2.Range("E11").Value = 4.7
$ws2.Range("F11").Value = 6

$ws2.Range("A12").Value = 8
$ws2.Range("B12").Value = "bottom component"
$ws2.Range("F12").Value = 1
$ws2.Range("H12").Value = 8
$ws2.Range("I12").Value = 54

# ---- total ------------------------------------------------------------
$ws2.Range("F19").Formula = "=+SUM(F2:F12)"

# ---- column widths (matches the "8-layer" sheet's best-fit look) ------
$ws2.Columns.Item(1).ColumnWidth = 3.5
$ws2.Columns.Item(2).ColumnWidth = 17.666666666666668
$ws2.Columns.Item(3).ColumnWidth = 13.833333333333334
$ws2.Columns.Item(4).ColumnWidth = 13.666666666666666
$ws2.Columns.Item(5).ColumnWidth = 3.1666666666666665
$ws2.Columns.Item(6).ColumnWidth = 8.666666666666666
$ws2.Columns.Item(7).ColumnWidth = 9.5
$ws2.Columns.Item(8).ColumnWidth = 10.666666666666666
$ws2.Columns.Item(9).ColumnWidth = 11.833333333333334

# ---- view state ---------------------------------------------------------
# "8-layer" loses the tab selection / single-cell selection, gaining a
# full used-range selection; "6-layer" becomes the active tab with E16
# selected.
$ws1.Range("A1:I23").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("E16").Select() | Out-Null
